# Generate Report for Handoff
# Moves the localization status from "In Translation" to "Ready for handoff"
# and refreshes the handoff timestamps on the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status text: "In Translation" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Latest handoff / HO xliff generation timestamps
$wsZhCn.Range("H2").Value = "2016-08-18 00:36:50"
$wsDeDe.Range("H2").Value = "2016-08-18 00:36:56"
$wsOverview.Range("G2").Value = "2016-08-18 00:36:56"

# The longer "Ready for handoff" label no longer fits the old column width,
# so the Status-related columns grow to accommodate it.
$wsOverview.Columns("E").ColumnWidth = 16.33
$wsOverview.Columns("F").ColumnWidth = 16.33
$wsZhCn.Columns("C").ColumnWidth = 16.33
$wsDeDe.Columns("C").ColumnWidth = 16.33
